# Fix retailer column (Digikey "H") part numbers that were losing their
# dash because Excel/LibreOffice auto-coerced them to plain numbers
# (e.g. 6242496 instead of 624-2496). Re-enter the affected cells as text
# so the dash is preserved, and move the active selection to H10 as part
# of a quick manual check of the fix ("add some tests").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sheet's inherited default column-width override so the sheet
# falls back to Excel's normal default column formatting.
$ws.Cells.ClearFormats()

# These three retailer references were stored as numbers, silently
# dropping the "XXX-XXXX" dash. Re-enter them as text with the dash.
$ws.Range("H2").Value = "624-2496"
$ws.Range("H3").Value = "788-2893"
$ws.Range("H5").Value = "670-8826"

# Leave the selection on H10, like the manual check that was used while
# verifying the fix.
[void]$ws.Range("H10").Select()
